$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.812.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.632.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.17"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.70"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.105.79"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.737.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.625.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.17"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.51"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "343.75"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +9.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000113"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "583.05"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.27"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.65"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.49"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.21%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "163.56"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.13"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.93"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0587"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.635"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0249"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0237"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.789"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.54%  "
